$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value that was bumped from
# 2023-09-01 (45170) to 2023-09-05 (45174) for every data row (2-89).
$ws.Range("C2:C89").Value = 45174
